$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 200.66667
$ws.Range("I8").Value = 40.8
$ws.Range("K8").Value = 122.4
$ws.Range("M8").Value = 16.60000000000001
$ws.Range("H12").Value = 1609.0714
$ws.Range("I12").Value = 2001.5714
$ws.Range("J12").Value = 1216.5714
$ws.Range("K12").Value = 2001.5714
$ws.Range("L12").Value = 1216.5714
$ws.Range("M12").Value = -1831.5714
$ws.Range("N12").Value = -1556.5714
$ws.Range("H19").Value = 5525
$ws.Range("I19").Value = 10250
$ws.Range("J19").Value = 800
$ws.Range("K19").Value = 10250
$ws.Range("L19").Value = 800
$ws.Range("M19").Value = -10075
$ws.Range("N19").Value = -1150
$ws.Range("H33").Value = 115.52381
$ws.Range("I33").Value = 120.210526
$ws.Range("K33").Value = 120.210526
$ws.Range("M33").Value = 108.789474
$ws.Range("H51").Value = 2967
$ws.Range("I51").Value = 3500.5
$ws.Range("J51").Value = 1900
$ws.Range("K51").Value = 3500.5
$ws.Range("L51").Value = 1900
$ws.Range("M51").Value = -3016.5
$ws.Range("N51").Value = -2868
$ws.Range("H70").Value = 908.3333
$ws.Range("J70").Value = 890
$ws.Range("L70").Value = 2670
$ws.Range("N70").Value = -3210
$ws.Range("H73").Value = 908.3333
$ws.Range("J73").Value = 890
$ws.Range("L73").Value = 2670
$ws.Range("N73").Value = -4542
$ws.Range("H88").Value = 2036.5
$ws.Range("I88").Value = 1983.3334
$ws.Range("K88").Value = 1983.3334
$ws.Range("M88").Value = -1577.3334
$ws.Range("H91").Value = 2036.5
$ws.Range("I91").Value = 1983.3334
$ws.Range("K91").Value = 1983.3334
$ws.Range("M91").Value = -579.3334
$ws.Range("H132").Value = 2045.7174
$ws.Range("I132").Value = 2046.7333
$ws.Range("K132").Value = 6140.199900000001
$ws.Range("M132").Value = -3610.199900000001
$ws.Range("H137").Value = 2769.476
$ws.Range("I137").Value = 2668.1765
$ws.Range("J137").Value = 3200
$ws.Range("K137").Value = 8004.529500000001
$ws.Range("L137").Value = 9600
$ws.Range("M137").Value = -5454.529500000001
$ws.Range("N137").Value = -14700

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2289.28
$ws.Range("I61").Value = 1904.5
$ws.Range("J61").Value = 5111
$ws.Range("K61").Value = 1904.5
$ws.Range("L61").Value = 5111
$ws.Range("M61").Value = -1692.5
$ws.Range("N61").Value = -5535
$ws.Range("H63").Value = 2124.5
$ws.Range("I63").Value = 1999.3334
$ws.Range("K63").Value = 1999.3334
$ws.Range("M63").Value = -1313.3334
$ws.Range("H66").Value = 2124.5
$ws.Range("I66").Value = 1999.3334
$ws.Range("K66").Value = 9996.666999999999
$ws.Range("M66").Value = -6564.666999999999
$ws.Range("H74").Value = 55556760
$ws.Range("I74").Value = 76923784
$ws.Range("K74").Value = 76923784
$ws.Range("M74").Value = -76922910
$ws.Range("H77").Value = 55556760
$ws.Range("I77").Value = 76923784
$ws.Range("K77").Value = 384618920
$ws.Range("M77").Value = -384614552
$ws.Range("H132").Value = 13099.863
$ws.Range("I132").Value = 1682.375
$ws.Range("K132").Value = 5047.125
$ws.Range("M132").Value = -2517.125
$ws.Range("H136").Value = 2289.28
$ws.Range("I136").Value = 1904.5
$ws.Range("J136").Value = 5111
$ws.Range("K136").Value = 5713.5
$ws.Range("L136").Value = 15333
$ws.Range("M136").Value = -3163.5
$ws.Range("N136").Value = -20433

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 30260
$ws.Range("J126").Value = 30260
$ws.Range("L126").Value = 30260
$ws.Range("N126").Value = -40140
$ws.Range("H134").Value = 3964.724
$ws.Range("I134").Value = 4419.08
$ws.Range("K134").Value = 13257.24
$ws.Range("M134").Value = -10722.24

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 882.3333
$ws.Range("I16").Value = 848.7143
$ws.Range("K16").Value = 848.7143
$ws.Range("M16").Value = -561.7143
$ws.Range("H58").Value = 14943.473
$ws.Range("I58").Value = 1086.6666
$ws.Range("J58").Value = 42657.082
$ws.Range("K58").Value = 1086.6666
$ws.Range("L58").Value = 42657.082
$ws.Range("M58").Value = -883.6666
$ws.Range("N58").Value = -43063.082
$ws.Range("H113").Value = 882.3333
$ws.Range("I113").Value = 848.7143
$ws.Range("K113").Value = 848.7143
$ws.Range("M113").Value = 1321.2857
$ws.Range("H132").Value = 11495.34
$ws.Range("I132").Value = 14672.737
$ws.Range("K132").Value = 44018.211
$ws.Range("M132").Value = -41488.211
$ws.Range("H134").Value = 1089.0938
$ws.Range("I134").Value = 828.9756
$ws.Range("K134").Value = 2486.9268
$ws.Range("M134").Value = 48.07319999999982
$ws.Range("H136").Value = 14943.473
$ws.Range("I136").Value = 1086.6666
$ws.Range("J136").Value = 42657.082
$ws.Range("K136").Value = 3259.9998
$ws.Range("L136").Value = 127971.246
$ws.Range("M136").Value = -709.9998000000001
$ws.Range("N136").Value = -133071.246

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 120.545456
$ws.Range("I12").Value = 50
$ws.Range("J12").Value = 136.22223
$ws.Range("K12").Value = 150
$ws.Range("L12").Value = 408.66669
$ws.Range("M12").Value = 23
$ws.Range("N12").Value = -754.66669
$ws.Range("H36").Value = 173058.14
$ws.Range("I36").Value = 3802
$ws.Range("J36").Value = 240760.6
$ws.Range("K36").Value = 11406
$ws.Range("L36").Value = 722281.8
$ws.Range("M36").Value = -11237
$ws.Range("N36").Value = -722619.8
$ws.Range("H131").Value = 769.41
$ws.Range("J131").Value = 769.41
$ws.Range("L131").Value = 2308.23
$ws.Range("N131").Value = -12388.23

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 19790.742
$ws.Range("J132").Value = 66082.25
$ws.Range("L132").Value = 198246.75
$ws.Range("N132").Value = -203306.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 732.5833
$ws.Range("J46").Value = 925
$ws.Range("L46").Value = 925
$ws.Range("N46").Value = -1301
$ws.Range("H61").Value = 3876.8076
$ws.Range("I61").Value = 2100.2
$ws.Range("J61").Value = 9798.833000000001
$ws.Range("K61").Value = 2100.2
$ws.Range("L61").Value = 9798.833000000001
$ws.Range("M61").Value = -1898.2
$ws.Range("N61").Value = -10202.833
$ws.Range("H68").Value = 2378.4
$ws.Range("I68").Value = 2214.6667
$ws.Range("J68").Value = 2624
$ws.Range("K68").Value = 2214.6667
$ws.Range("L68").Value = 2624
$ws.Range("M68").Value = -1465.6667
$ws.Range("N68").Value = -4122
$ws.Range("H71").Value = 2378.4
$ws.Range("I71").Value = 2214.6667
$ws.Range("J71").Value = 2624
$ws.Range("K71").Value = 11073.3335
$ws.Range("L71").Value = 13120
$ws.Range("M71").Value = -7329.333500000001
$ws.Range("N71").Value = -20608
$ws.Range("H93").Value = 2263.3076
$ws.Range("I93").Value = 2454.3
$ws.Range("J93").Value = 1626.6666
$ws.Range("K93").Value = 2454.3
$ws.Range("L93").Value = 1626.6666
$ws.Range("M93").Value = -1206.3
$ws.Range("N93").Value = -4122.6666
$ws.Range("H113").Value = 3876.8076
$ws.Range("I113").Value = 2100.2
$ws.Range("J113").Value = 9798.833000000001
$ws.Range("K113").Value = 2100.2
$ws.Range("L113").Value = 9798.833000000001
$ws.Range("M113").Value = 69.80000000000018
$ws.Range("N113").Value = -14138.833
$ws.Range("H122").Value = 1156073.4
$ws.Range("I122").Value = 1785049.9
$ws.Range("J122").Value = 2950
$ws.Range("K122").Value = 5355149.699999999
$ws.Range("L122").Value = 8850
$ws.Range("M122").Value = -5352699.699999999
$ws.Range("N122").Value = -13750
$ws.Range("H132").Value = 2328.7896
$ws.Range("J132").Value = 3999.875
$ws.Range("L132").Value = 11999.625
$ws.Range("N132").Value = -17059.625
$ws.Range("H136").Value = 28805.277
$ws.Range("I136").Value = 34380
$ws.Range("J136").Value = 931.6667
$ws.Range("K136").Value = 103140
$ws.Range("L136").Value = 2795.0001
$ws.Range("M136").Value = -100590
$ws.Range("N136").Value = -7895.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 4546170.5
$ws.Range("I107").Value = 1040
$ws.Range("K107").Value = 3120
$ws.Range("M107").Value = -1200
$ws.Range("H113").Value = 1423247
$ws.Range("I113").Value = 1113.6364
$ws.Range("J113").Value = 3378680.5
$ws.Range("K113").Value = 3340.9092
$ws.Range("L113").Value = 10136041.5
$ws.Range("M113").Value = -1170.9092
$ws.Range("N113").Value = -10140381.5
$ws.Range("H122").Value = 1789
$ws.Range("I122").Value = 1816.4286
$ws.Range("J122").Value = 1717
$ws.Range("K122").Value = 5449.2858
$ws.Range("L122").Value = 5151
$ws.Range("M122").Value = -2999.2858
$ws.Range("N122").Value = -10051
$ws.Range("H132").Value = 1486.1305
$ws.Range("I132").Value = 1131.3529
$ws.Range("J132").Value = 2491.3333
$ws.Range("K132").Value = 3394.0587
$ws.Range("L132").Value = 7473.999899999999
$ws.Range("M132").Value = -864.0587000000005
$ws.Range("N132").Value = -12533.9999
$ws.Range("H136").Value = 38463810
$ws.Range("I136").Value = 66668856
$ws.Range("K136").Value = 200006568
$ws.Range("M136").Value = -200004018
